$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.237799999999997
$ws.Range("C3").Value = -10.24699999999999
$ws.Range("D3").Value = -7.117499999999993
$ws.Range("C4").Value = -12.5522
$ws.Range("D9").Value = -7.099199999999997
$ws.Range("A11").Value = -21.84850000000001
$ws.Range("A12").Value = -21.5158
$ws.Range("C14").Value = -13.2525
$ws.Range("A15").Value = -21.78399999999999
$ws.Range("D15").Value = -8.149099999999997
$ws.Range("D19").Value = -7.635399999999995
$ws.Range("D20").Value = -7.566599999999999
$ws.Range("D25").Value = -7.596400000000002
$ws.Range("C26").Value = -12.88400000000001
$ws.Range("A27").Value = -21.91160000000001
$ws.Range("D27").Value = -8.757599999999998
$ws.Range("A28").Value = -21.73409999999999
$ws.Range("D28").Value = -7.980399999999999
$ws.Range("D30").Value = -7.066300000000004
$ws.Range("A31").Value = -21.539
$ws.Range("C31").Value = -13.021
$ws.Range("A32").Value = -21.66300000000001
$ws.Range("D32").Value = -8.253500000000006
$ws.Range("C35").Value = -12.52080000000001
$ws.Range("A36").Value = -20.6175
$ws.Range("C37").Value = -13.86779999999999
$ws.Range("A38").Value = -19.3145
$ws.Range("C39").Value = -12.87780000000001
$ws.Range("C40").Value = -14.0709
$ws.Range("D44").Value = -7.258700000000002
$ws.Range("C45").Value = -14.29609999999999
$ws.Range("A46").Value = -21.59000000000001
$ws.Range("D47").Value = -7.481699999999997
$ws.Range("C52").Value = -10.8627
$ws.Range("A54").Value = -21.5511
$ws.Range("A55").Value = -22.19670000000001
$ws.Range("A56").Value = -22.10990000000001
$ws.Range("C57").Value = -14.59829999999999
$ws.Range("D58").Value = -8.037599999999998
$ws.Range("D62").Value = -8.290499999999987
$ws.Range("A67").Value = -21.52909999999997
$ws.Range("A69").Value = -21.75099999999998
$ws.Range("A72").Value = -21.48089999999999
$ws.Range("A73").Value = -19.96840000000001
$ws.Range("D77").Value = -5.7643
$ws.Range("D78").Value = -7.451700000000001
$ws.Range("C81").Value = -12.7172
$ws.Range("A83").Value = -21.65309999999999
$ws.Range("C83").Value = -12.01090000000001
$ws.Range("D84").Value = -8.601099999999997
$ws.Range("A86").Value = -22.15780000000001
$ws.Range("D89").Value = -6.129499999999998
$ws.Range("A91").Value = -21.4587
$ws.Range("D91").Value = -6.126999999999996
$ws.Range("D92").Value = -6.015199999999997
$ws.Range("A93").Value = -21.29419999999999
$ws.Range("D96").Value = -7.377700000000004
$ws.Range("A99").Value = -20.52209999999998
$ws.Range("C100").Value = -12.5552
$ws.Range("C102").Value = -14.375
$ws.Range("D102").Value = -7.889700000000001
